$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "genre" column before the existing "watched" column (D),
# which shifts watched -> E, rate -> F, review -> G.
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "genre"

# Row 2 (Interstellar) updates: year corrected, rating lowered, review text added.
# Keep the year cell text-typed (not numeric) to match the original data type.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2018"
$ws.Range("C2").Style = "Normal"

# D2 (genre) intentionally left blank/unfinished per the commit message.

$ws.Range("F2").Value = "★★★★★ (5.0)"
$ws.Range("G2").Value = "Very good! i like how cooper eventually meet his daughter again"

# New row 3: Taxi Driver
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Taxi Driver"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "199x"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "Action, Loneliness"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "-"
$ws.Range("G3").Value = "-"
